# #5: cash & deposit done
# Turns the "存款" (deposit) sheet's ad-hoc 6-column dump into the same
# normalized schema used by the other property sheets: adds a proper
# header row and appends property_category/category/date/legislator_name/
# legislator_id/source_file/index columns (G:M) to every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("存款")

# ---- Header row (row 1): replace the stray duplicated data with column
# names, and extend it through column M -------------------------------
$ws.Cells.Item(1, 2).Value = "bank"
$ws.Cells.Item(1, 3).Value = "deposit_type"
$ws.Cells.Item(1, 4).Value = "currency"
$ws.Cells.Item(1, 5).Value = "owner"
$ws.Cells.Item(1, 6).Value = "total"
$ws.Cells.Item(1, 7).Value = "property_category"
$ws.Cells.Item(1, 8).Value = "category"
$ws.Cells.Item(1, 9).Value = "date"
$ws.Cells.Item(1, 10).Value = "legislator_name"
$ws.Cells.Item(1, 11).Value = "legislator_id"
$ws.Cells.Item(1, 12).Value = "source_file"
$ws.Cells.Item(1, 13).Value = "index"

# Give the new header cells (G1:M1) the same bold / centered / boxed look
# as the existing header cells (B1:F1, style index 1).
$headerRng = $ws.Range("G1:M1")
$headerRng.Font.Bold = $true
$headerRng.HorizontalAlignment = -4108
$headerRng.VerticalAlignment = -4160
$headerRng.Borders.LineStyle = 1

# ---- Data rows (2-5): fix B so it doesn't overwrite with row1's bank
# name, make F numeric everywhere, and append the common metadata
# columns G:M ------------------------------------------------------------
$banks = @("臺灣銀行群賢分行", "臺灣銀行汐止分行", "永豐商業銀行中正分行", "台北富邦商業銀行復興分行")
$amounts = @(7032144, 2648902, 97957, 623009)

for ($i = 0; $i -lt 4; $i++) {
    $r = $i + 2

    $ws.Cells.Item($r, 2).Value = $banks[$i]
    $ws.Cells.Item($r, 3).Value = "活期存款"
    $ws.Cells.Item($r, 4).Value = "新臺幣"
    $ws.Cells.Item($r, 5).Value = "李慶華"
    $ws.Cells.Item($r, 6).Value = $amounts[$i]

    $ws.Cells.Item($r, 7).Value = "deposit"
    $ws.Cells.Item($r, 8).Value = "normal"
    $ws.Cells.Item($r, 9).Value = "2011-11-30"
    $ws.Cells.Item($r, 10).Value = "李慶華"
    $ws.Cells.Item($r, 11).Value = 607
    $ws.Cells.Item($r, 12).Value = "tmp68f81"
    $ws.Cells.Item($r, 13).Value = $r + 56
}
